$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': 10, 'learning_rate': 'constant', 'learning_rate_init': 0.001, 'max_iter': 1900, 'solver': 'lbfgs'}"
$ws.Range("B9").Value = 0.975
$ws.Range("C9").Value = 0.783
